$d = $word.ActiveDocument

# 1. Merge "{{ " / "stopien" / " }} {{ " / "imie" / " }} {{ nazwisko }}" runs
#    (with spell-check proofErr markers) into a single run's text.
$d.Content.Find.Execute("{{ stopien }} {{ imie }} {{ nazwisko }}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{{ stopien }} {{ imie }} {{ nazwisko }}", 2) | Out-Null

# 2. Merge "{{ pluton }}" / "pl" / "/5kmp/2BS" runs into a single run's text.
$d.Content.Find.Execute("{{ pluton }}pl/5kmp/2BS", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{{ pluton }}pl/5kmp/2BS", 2) | Out-Null

# 3. Merge "{{ " / "data_przepustki" / " }}" runs into a single run's text.
$d.Content.Find.Execute("{{ data_przepustki }}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{{ data_przepustki }}", 2) | Out-Null

# 4. Merge "{{ " / "zaleglosci" / " }}" runs into a single run's text.
$d.Content.Find.Execute("{{ zaleglosci }}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{{ zaleglosci }}", 2) | Out-Null

# 5. Fix "{{ miejscowość }}." -> "{{ miejscowosc }}." (drop Polish diacritics),
#    splitting the run into three pieces as in the target edit.
$d.Content.Find.Execute("{{ miejscowość }}.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{{ miejscowosc }}.", 2) | Out-Null
